$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "Closed"
$ws.Range("H3").Value = "Closed"
$ws.Range("H4").Value = "Closed"

$ws.Range("H4").Select() | Out-Null
